$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.933.56"
$ws.Range("E2").Value = "  +0.35%  "
$ws.Range("D3").Value = "1.905.86"
$ws.Range("E3").Value = "  +0.65%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7977"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.18%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "241.68"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.07%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3154"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.44%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "26.22"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.65%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06904"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.13%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07989"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.17%  "
$ws.Range("D12").Value = "1.901.81"
$ws.Range("E12").Value = "  +0.52%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.7387"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.44%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.181"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.33%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "92.78"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.92%  "
$ws.Range("D16").Value = "29.923.32"
$ws.Range("E16").Value = "  +0.35%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.93"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.40%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.848"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.81%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "244.87"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.73%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007725"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.66%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.000"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.02%  "
$ws.Range("D22").Value = "2.148.99"
$ws.Range("E22").Value = "  +0.22%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.0000"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.805"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.85%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "167.65"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.48%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.182"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.73%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1403"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +9.00%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.86"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.05%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.021"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.90%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.362"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.71%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.512"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.02%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.295"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.37%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.078"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.40%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05461"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.14%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.254"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.88%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7264"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.41%  "
$ws.Range("E37").Value = "  +0.21%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01922"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.31%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.781"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.51%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.133"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.12%  "
$ws.Range("E41").Value = "  +0.03%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "71.95"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.58%  "
$ws.Range("E43").Value = "  -0.04%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8345"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.33%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.865"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.60%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "100.19"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.83%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.507"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.95%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.688"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.95%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "986.72"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +7.90%  "
$ws.Range("D50").Value = "2.052.69"
$ws.Range("E50").Value = "  +0.13%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "36.09"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.33%  "
